$wb = $excel.ActiveWorkbook
$origActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("SearchHomePage")

$ws.Range("A4").Value = "Sanity_12"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "Accleration"

$ws.Range("A5").Value = "Sanity_13"
$ws.Range("B5").Value = "Yes"

$ws.Range("C5").Select()

$origActive.Activate()
